$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# 1. Fix the "Path Materialization" bullet text on the left content placeholder.
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$old = "Path Materialization: if exploring a virtual edge, materialize its virtual nodes;"
$new = "Path instantiation: if exploring a virtual edge, instantiation its virtual nodes;"
$full = $tr.Text
$idx = $full.IndexOf($old)
if ($idx -ge 0) {
    $c = $tr.Characters($idx + 1, $old.Length)
    $c.Text = $new
}

# 2. Remove the slide's click-animation timing (the <p:timing> build sequence)
#    by deleting every effect from the slide's main animation sequence.
$seq = $s.TimeLine.MainSequence
while ($seq.Count -gt 0) {
    $seq.Item(1).Delete()
}
